# Scheduled runner update: refresh market-price derived columns (H-N)
# for FFXIV leve-profit tracking sheets (currentAveragePrice*, LevePrice*, LeveProfit*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1390209.9
$ws.Range("I103").Value = 2778719.8
$ws.Range("J103").Value = 1700
$ws.Range("K103").Value = 8336159.399999999
$ws.Range("L103").Value = 5100
$ws.Range("M103").Value = -8335573.399999999
$ws.Range("N103").Value = -6272
$ws.Range("H114").Value = 36321
$ws.Range("J114").Value = 36321
$ws.Range("L114").Value = 36321
$ws.Range("N114").Value = -44999
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H127").Value = 68152.92999999999
$ws.Range("I127").Value = 695.4
$ws.Range("J127").Value = 101881.7
$ws.Range("K127").Value = 2086.2
$ws.Range("L127").Value = 305645.1
$ws.Range("M127").Value = 2873.8
$ws.Range("N127").Value = -315565.1
$ws.Range("H132").Value = 1318644.1
$ws.Range("I132").Value = 1541556.9
$ws.Range("J132").Value = 3459
$ws.Range("K132").Value = 4624670.699999999
$ws.Range("L132").Value = 10377
$ws.Range("M132").Value = -4622140.699999999
$ws.Range("N132").Value = -15437

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1174.1562
$ws.Range("I61").Value = 951.73914
$ws.Range("J61").Value = 1742.5555
$ws.Range("K61").Value = 951.73914
$ws.Range("L61").Value = 1742.5555
$ws.Range("M61").Value = -739.73914
$ws.Range("N61").Value = -2166.5555
$ws.Range("H74").Value = 10639078
$ws.Range("I74").Value = 12195827
$ws.Range("J74").Value = 1300
$ws.Range("K74").Value = 12195827
$ws.Range("L74").Value = 1300
$ws.Range("M74").Value = -12194953
$ws.Range("N74").Value = -3048
$ws.Range("H77").Value = 10639078
$ws.Range("I77").Value = 12195827
$ws.Range("J77").Value = 1300
$ws.Range("K77").Value = 60979135
$ws.Range("L77").Value = 6500
$ws.Range("M77").Value = -60974767
$ws.Range("N77").Value = -15236
$ws.Range("H102").Value = 1495
$ws.Range("I102").Value = 1495
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1495
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 127
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 1603.0435
$ws.Range("I132").Value = 925.1515000000001
$ws.Range("J132").Value = 3323.8462
$ws.Range("K132").Value = 2775.4545
$ws.Range("L132").Value = 9971.5386
$ws.Range("M132").Value = -245.4545000000003
$ws.Range("N132").Value = -15031.5386
$ws.Range("H136").Value = 1174.1562
$ws.Range("I136").Value = 951.73914
$ws.Range("J136").Value = 1742.5555
$ws.Range("K136").Value = 2855.21742
$ws.Range("L136").Value = 5227.666499999999
$ws.Range("M136").Value = -305.2174199999999
$ws.Range("N136").Value = -10327.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 25646860
$ws.Range("I20").Value = 43485588
$ws.Range("J20").Value = 3692.875
$ws.Range("K20").Value = 43485588
$ws.Range("L20").Value = 3692.875
$ws.Range("M20").Value = -43485341
$ws.Range("N20").Value = -4186.875
$ws.Range("H105").Value = 875396.5600000001
$ws.Range("I105").Value = 1083424.4
$ws.Range("J105").Value = 1679.8
$ws.Range("K105").Value = 1083424.4
$ws.Range("L105").Value = 1679.8
$ws.Range("M105").Value = -1081677.4
$ws.Range("N105").Value = -5173.8
$ws.Range("H134").Value = 8929474
$ws.Range("I134").Value = 9260184
$ws.Range("J134").Value = 314
$ws.Range("K134").Value = 27780552
$ws.Range("L134").Value = 942
$ws.Range("M134").Value = -27778017
$ws.Range("N134").Value = -6012

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1770.5853
$ws.Range("I132").Value = 1277.4482
$ws.Range("J132").Value = 2962.3333
$ws.Range("K132").Value = 3832.3446
$ws.Range("L132").Value = 8886.999899999999
$ws.Range("M132").Value = -1302.3446
$ws.Range("N132").Value = -13946.9999
$ws.Range("H134").Value = 1792.0312
$ws.Range("I134").Value = 1486.381
$ws.Range("J134").Value = 2375.5454
$ws.Range("K134").Value = 4459.143
$ws.Range("L134").Value = 7126.6362
$ws.Range("M134").Value = -1924.143
$ws.Range("N134").Value = -12196.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 1090.625
$ws.Range("J104").Value = 1124.1666
$ws.Range("L104").Value = 3372.4998
$ws.Range("N104").Value = -8614.4998
$ws.Range("H131").Value = 727.9114
$ws.Range("I131").Value = 284.72
$ws.Range("J131").Value = 933.0925999999999
$ws.Range("K131").Value = 854.1600000000001
$ws.Range("L131").Value = 2799.2778
$ws.Range("M131").Value = 4185.84
$ws.Range("N131").Value = -12879.2778

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1590.6086
$ws.Range("I122").Value = 1505.0667
$ws.Range("J122").Value = 1751
$ws.Range("K122").Value = 4515.2001
$ws.Range("L122").Value = 5253
$ws.Range("M122").Value = -2065.2001
$ws.Range("N122").Value = -10153

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1637.125
$ws.Range("I7").Value = 1514.9231
$ws.Range("J7").Value = 2166.6667
$ws.Range("K7").Value = 1514.9231
$ws.Range("L7").Value = 2166.6667
$ws.Range("M7").Value = -1402.9231
$ws.Range("N7").Value = -2390.6667
$ws.Range("H22").Value = 322.82053
$ws.Range("I22").Value = 282.66666
$ws.Range("J22").Value = 456.66666
$ws.Range("K22").Value = 282.66666
$ws.Range("L22").Value = 456.66666
$ws.Range("M22").Value = 12.33334000000002
$ws.Range("N22").Value = -1046.66666
$ws.Range("H27").Value = 322.82053
$ws.Range("I27").Value = 282.66666
$ws.Range("J27").Value = 456.66666
$ws.Range("K27").Value = 282.66666
$ws.Range("L27").Value = 456.66666
$ws.Range("M27").Value = -175.66666
$ws.Range("N27").Value = -670.66666
$ws.Range("H100").Value = 2082.0908
$ws.Range("I100").Value = 1460.6
$ws.Range("J100").Value = 2600
$ws.Range("K100").Value = 1460.6
$ws.Range("L100").Value = 2600
$ws.Range("M100").Value = -919.5999999999999
$ws.Range("N100").Value = -3682
$ws.Range("H122").Value = 4928.234
$ws.Range("I122").Value = 5120.049
$ws.Range("J122").Value = 3617.5
$ws.Range("K122").Value = 15360.147
$ws.Range("L122").Value = 10852.5
$ws.Range("M122").Value = -12910.147
$ws.Range("N122").Value = -15752.5
$ws.Range("H126").Value = 1637.125
$ws.Range("I126").Value = 1514.9231
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 4544.7693
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -2074.7693
$ws.Range("N126").Value = -11440.0001
$ws.Range("H136").Value = 2157.4807
$ws.Range("I136").Value = 1654.5428
$ws.Range("J136").Value = 3192.9412
$ws.Range("K136").Value = 4963.6284
$ws.Range("L136").Value = 9578.8236
$ws.Range("M136").Value = -2413.6284
$ws.Range("N136").Value = -14678.8236

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 689.8395
$ws.Range("I132").Value = 358.5098
$ws.Range("J132").Value = 1253.1
$ws.Range("K132").Value = 1075.5294
$ws.Range("L132").Value = 3759.3
$ws.Range("M132").Value = 1454.4706
$ws.Range("N132").Value = -8819.299999999999
